$wb = $excel.ActiveWorkbook

# ---- Sheet: Overview ----
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A2").Value = 'ffffb6ad0d8b-c1c0-4c1f-bfc6-69772c211997.md'
$ws.Range("B2").Value = 'Handed back: in sync with en-US'
$ws.Range("C2").Value = 'Handed back: in sync with en-US'
$ws.Range("D2").Value = '2016-02-13 05:02:08'
$ws.Range("A3").Value = 'ffffff434d3f9c-c5e3-4796-b805-00c9b17a6d93.md'
$ws.Range("B3").Value = 'Handed back: in sync with en-US'
$ws.Range("C3").Value = 'Handed back: in sync with en-US'
$ws.Range("D3").Value = '2016-02-13 05:02:08'
$ws.Range("A4").Value = '7a01dfb2-5efe-44c9-9ceb-2d775ef4025f.md'
$ws.Range("B4").Value = 'Ready for handoff'
$ws.Range("C4").Value = 'Ready for handoff'
$ws.Range("D4").Value = '2016-04-13 05:04:03'

# Rebuild hyperlinks: delete all, then re-add in original order/URLs so
# relationship ids (rId2, rId3, ...) are reassigned identically; only the
# display text (matching the new cell content) changes.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), 'https://github.com/OpenLocalizationTest/oltest/blob/fe32e375100bc6b15da6c00d5bd65ac940e27279/e2e/7a01dfb2-5efe-44c9-9ceb-2d775ef4025f.md', [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, 'ffffb6ad0d8b-c1c0-4c1f-bfc6-69772c211997.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), 'https://github.com/OpenLocalizationTest/oltest/blob/fe32e375100bc6b15da6c00d5bd65ac940e27279/e2e/ffffb6ad0d8b-c1c0-4c1f-bfc6-69772c211997.md', [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, 'ffffff434d3f9c-c5e3-4796-b805-00c9b17a6d93.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), 'https://github.com/OpenLocalizationTest/oltest/blob/fe32e375100bc6b15da6c00d5bd65ac940e27279/e2e/ffffff434d3f9c-c5e3-4796-b805-00c9b17a6d93.md', [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, '7a01dfb2-5efe-44c9-9ceb-2d775ef4025f.md') | Out-Null

# ---- Sheet: zh-cn ----
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A2").Value = 'ffffb6ad0d8b-c1c0-4c1f-bfc6-69772c211997.md'
$ws.Range("B2").Value = '.md'
$ws.Range("C2").Value = 'Handed back: in sync with en-US'
$ws.Range("D2").Value = 'a3573822-9d94-4592-834d-ffb5b55a6ca7.2924ee7c1e8c42dfa5b48b47664a134ed1cea41e.zh-cn.xlf'
$ws.Range("E2").Value = '2016-03-13 05:02:04'
$ws.Range("F2").Value = 'a3573822-9d94-4592-834d-ffb5b55a6ca7.md'
$ws.Range("G2").Value = 'a3573822-9d94-4592-834d-ffb5b55a6ca7.2924ee7c1e8c42dfa5b48b47664a134ed1cea41e.zh-cn.xlf'
$ws.Range("H2").Value = '2016-03-13 05:02:22'
$ws.Range("I2").Value = 'Include'
$ws.Range("A3").Value = 'ffffff434d3f9c-c5e3-4796-b805-00c9b17a6d93.md'
$ws.Range("B3").Value = '.md'
$ws.Range("C3").Value = 'Handed back: in sync with en-US'
$ws.Range("D3").Value = 'a3573822-9d94-4592-834d-ffb5b55a6ca7.2924ee7c1e8c42dfa5b48b47664a134ed1cea41e.zh-cn.xlf'
$ws.Range("E3").Value = '2016-03-13 05:02:04'
$ws.Range("F3").Value = 'a3573822-9d94-4592-834d-ffb5b55a6ca7.md'
$ws.Range("G3").Value = 'a3573822-9d94-4592-834d-ffb5b55a6ca7.2924ee7c1e8c42dfa5b48b47664a134ed1cea41e.zh-cn.xlf'
$ws.Range("H3").Value = '2016-03-13 05:02:22'
$ws.Range("I3").Value = 'Include'
$ws.Range("A4").Value = '7a01dfb2-5efe-44c9-9ceb-2d775ef4025f.md'
$ws.Range("B4").Value = '.md'
$ws.Range("C4").Value = 'Ready for handoff'
$ws.Range("D4").Value = '7a01dfb2-5efe-44c9-9ceb-2d775ef4025f.718e13acabd4baa934d39fc555b8ee9c989bf7e0.zh-cn.xlf'
$ws.Range("E4").Value = '2016-03-13 05:03:56'
$ws.Range("F4").Value = '7a01dfb2-5efe-44c9-9ceb-2d775ef4025f.md'
$ws.Range("G4").Value = '7a01dfb2-5efe-44c9-9ceb-2d775ef4025f.718e13acabd4baa934d39fc555b8ee9c989bf7e0.zh-cn.xlf'
$ws.Range("H4").Value = '2016-03-13 05:03:23'
$ws.Range("I4").Value = 'Include'

# Rebuild hyperlinks: delete all, then re-add in original order/URLs so
# relationship ids (rId2, rId3, ...) are reassigned identically; only the
# display text (matching the new cell content) changes.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), 'https://github.com/OpenLocalizationTest/oltest/blob/fe32e375100bc6b15da6c00d5bd65ac940e27279/e2e/7a01dfb2-5efe-44c9-9ceb-2d775ef4025f.md', [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, 'ffffb6ad0d8b-c1c0-4c1f-bfc6-69772c211997.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range("B2"), 'https://github.com/OpenLocalizationTest/oltest/blob/fe32e375100bc6b15da6c00d5bd65ac940e27279/e2e/7a01dfb2-5efe-44c9-9ceb-2d775ef4025f.md', [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, '.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range("D2"), 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c39e7dd0a475e0a80fe24dbc68c73590408b4825/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/7a01dfb2-5efe-44c9-9ceb-2d775ef4025f.718e13acabd4baa934d39fc555b8ee9c989bf7e0.zh-cn.xlf', [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, 'a3573822-9d94-4592-834d-ffb5b55a6ca7.2924ee7c1e8c42dfa5b48b47664a134ed1cea41e.zh-cn.xlf') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F2"), 'https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/8a37dfec9c577dcd6e222db4144d1deb4b64868c/e2e/7a01dfb2-5efe-44c9-9ceb-2d775ef4025f.md', [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, 'a3573822-9d94-4592-834d-ffb5b55a6ca7.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range("G2"), 'https://github.com/OpenLocalizationTestOrg/olhandback/blob/38022ba51b03e529f96033ab250af84a0a0bcc83/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/7a01dfb2-5efe-44c9-9ceb-2d775ef4025f.718e13acabd4baa934d39fc555b8ee9c989bf7e0.zh-cn.xlf', [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, 'a3573822-9d94-4592-834d-ffb5b55a6ca7.2924ee7c1e8c42dfa5b48b47664a134ed1cea41e.zh-cn.xlf') | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), 'https://github.com/OpenLocalizationTest/oltest/blob/fe32e375100bc6b15da6c00d5bd65ac940e27279/e2e/ffffb6ad0d8b-c1c0-4c1f-bfc6-69772c211997.md', [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, 'ffffff434d3f9c-c5e3-4796-b805-00c9b17a6d93.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range("B3"), 'https://github.com/OpenLocalizationTest/oltest/blob/fe32e375100bc6b15da6c00d5bd65ac940e27279/e2e/ffffb6ad0d8b-c1c0-4c1f-bfc6-69772c211997.md', [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, '.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range("D3"), 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a57e240e148a8297ef3beed80705e24358c78f06/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a3573822-9d94-4592-834d-ffb5b55a6ca7.2924ee7c1e8c42dfa5b48b47664a134ed1cea41e.zh-cn.xlf', [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, 'a3573822-9d94-4592-834d-ffb5b55a6ca7.2924ee7c1e8c42dfa5b48b47664a134ed1cea41e.zh-cn.xlf') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F3"), 'https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/e5927f08c8836525f5f7e1dae0ac9c2fdb3d2682/e2e/a3573822-9d94-4592-834d-ffb5b55a6ca7.md', [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, 'a3573822-9d94-4592-834d-ffb5b55a6ca7.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range("G3"), 'https://github.com/OpenLocalizationTestOrg/olhandback/blob/f238b79cad5604da4013a49791a527474359b330/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a3573822-9d94-4592-834d-ffb5b55a6ca7.2924ee7c1e8c42dfa5b48b47664a134ed1cea41e.zh-cn.xlf', [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, 'a3573822-9d94-4592-834d-ffb5b55a6ca7.2924ee7c1e8c42dfa5b48b47664a134ed1cea41e.zh-cn.xlf') | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), 'https://github.com/OpenLocalizationTest/oltest/blob/fe32e375100bc6b15da6c00d5bd65ac940e27279/e2e/ffffff434d3f9c-c5e3-4796-b805-00c9b17a6d93.md', [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, '7a01dfb2-5efe-44c9-9ceb-2d775ef4025f.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range("B4"), 'https://github.com/OpenLocalizationTest/oltest/blob/fe32e375100bc6b15da6c00d5bd65ac940e27279/e2e/ffffff434d3f9c-c5e3-4796-b805-00c9b17a6d93.md', [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, '.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range("D4"), 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a57e240e148a8297ef3beed80705e24358c78f06/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a3573822-9d94-4592-834d-ffb5b55a6ca7.2924ee7c1e8c42dfa5b48b47664a134ed1cea41e.zh-cn.xlf', [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, '7a01dfb2-5efe-44c9-9ceb-2d775ef4025f.718e13acabd4baa934d39fc555b8ee9c989bf7e0.zh-cn.xlf') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F4"), 'https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/e5927f08c8836525f5f7e1dae0ac9c2fdb3d2682/e2e/a3573822-9d94-4592-834d-ffb5b55a6ca7.md', [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, '7a01dfb2-5efe-44c9-9ceb-2d775ef4025f.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range("G4"), 'https://github.com/OpenLocalizationTestOrg/olhandback/blob/f238b79cad5604da4013a49791a527474359b330/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a3573822-9d94-4592-834d-ffb5b55a6ca7.2924ee7c1e8c42dfa5b48b47664a134ed1cea41e.zh-cn.xlf', [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, '7a01dfb2-5efe-44c9-9ceb-2d775ef4025f.718e13acabd4baa934d39fc555b8ee9c989bf7e0.zh-cn.xlf') | Out-Null

# ---- Sheet: de-de ----
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A2").Value = 'ffffb6ad0d8b-c1c0-4c1f-bfc6-69772c211997.md'
$ws.Range("B2").Value = '.md'
$ws.Range("C2").Value = 'Handed back: in sync with en-US'
$ws.Range("D2").Value = 'a3573822-9d94-4592-834d-ffb5b55a6ca7.2924ee7c1e8c42dfa5b48b47664a134ed1cea41e.de-de.xlf'
$ws.Range("E2").Value = '2016-03-13 05:02:08'
$ws.Range("F2").Value = 'a3573822-9d94-4592-834d-ffb5b55a6ca7.md'
$ws.Range("G2").Value = 'a3573822-9d94-4592-834d-ffb5b55a6ca7.2924ee7c1e8c42dfa5b48b47664a134ed1cea41e.de-de.xlf'
$ws.Range("H2").Value = '2016-03-13 05:02:28'
$ws.Range("I2").Value = 'Include'
$ws.Range("A3").Value = 'ffffff434d3f9c-c5e3-4796-b805-00c9b17a6d93.md'
$ws.Range("B3").Value = '.md'
$ws.Range("C3").Value = 'Handed back: in sync with en-US'
$ws.Range("D3").Value = 'a3573822-9d94-4592-834d-ffb5b55a6ca7.2924ee7c1e8c42dfa5b48b47664a134ed1cea41e.de-de.xlf'
$ws.Range("E3").Value = '2016-03-13 05:02:08'
$ws.Range("F3").Value = 'a3573822-9d94-4592-834d-ffb5b55a6ca7.md'
$ws.Range("G3").Value = 'a3573822-9d94-4592-834d-ffb5b55a6ca7.2924ee7c1e8c42dfa5b48b47664a134ed1cea41e.de-de.xlf'
$ws.Range("H3").Value = '2016-03-13 05:02:28'
$ws.Range("I3").Value = 'Include'
$ws.Range("A4").Value = '7a01dfb2-5efe-44c9-9ceb-2d775ef4025f.md'
$ws.Range("B4").Value = '.md'
$ws.Range("C4").Value = 'Ready for handoff'
$ws.Range("D4").Value = '7a01dfb2-5efe-44c9-9ceb-2d775ef4025f.718e13acabd4baa934d39fc555b8ee9c989bf7e0.de-de.xlf'
$ws.Range("E4").Value = '2016-03-13 05:04:03'
$ws.Range("F4").Value = '7a01dfb2-5efe-44c9-9ceb-2d775ef4025f.md'
$ws.Range("G4").Value = '7a01dfb2-5efe-44c9-9ceb-2d775ef4025f.718e13acabd4baa934d39fc555b8ee9c989bf7e0.de-de.xlf'
$ws.Range("H4").Value = '2016-03-13 05:03:29'
$ws.Range("I4").Value = 'Include'

# Rebuild hyperlinks: delete all, then re-add in original order/URLs so
# relationship ids (rId2, rId3, ...) are reassigned identically; only the
# display text (matching the new cell content) changes.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), 'https://github.com/OpenLocalizationTest/oltest/blob/fe32e375100bc6b15da6c00d5bd65ac940e27279/e2e/7a01dfb2-5efe-44c9-9ceb-2d775ef4025f.md', [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, 'ffffb6ad0d8b-c1c0-4c1f-bfc6-69772c211997.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range("B2"), 'https://github.com/OpenLocalizationTest/oltest/blob/fe32e375100bc6b15da6c00d5bd65ac940e27279/e2e/7a01dfb2-5efe-44c9-9ceb-2d775ef4025f.md', [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, '.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range("D2"), 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/714b2537fdc7cc7a36de80bb174990ca67be79a3/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/7a01dfb2-5efe-44c9-9ceb-2d775ef4025f.718e13acabd4baa934d39fc555b8ee9c989bf7e0.de-de.xlf', [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, 'a3573822-9d94-4592-834d-ffb5b55a6ca7.2924ee7c1e8c42dfa5b48b47664a134ed1cea41e.de-de.xlf') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F2"), 'https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/387f2f572c957a39f38b2c46d6715d5559489d04/e2e/7a01dfb2-5efe-44c9-9ceb-2d775ef4025f.md', [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, 'a3573822-9d94-4592-834d-ffb5b55a6ca7.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range("G2"), 'https://github.com/OpenLocalizationTestOrg/olhandback/blob/b4aedd0204f41cebabbde194ffb31e07a2a41284/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/7a01dfb2-5efe-44c9-9ceb-2d775ef4025f.718e13acabd4baa934d39fc555b8ee9c989bf7e0.de-de.xlf', [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, 'a3573822-9d94-4592-834d-ffb5b55a6ca7.2924ee7c1e8c42dfa5b48b47664a134ed1cea41e.de-de.xlf') | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), 'https://github.com/OpenLocalizationTest/oltest/blob/fe32e375100bc6b15da6c00d5bd65ac940e27279/e2e/ffffb6ad0d8b-c1c0-4c1f-bfc6-69772c211997.md', [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, 'ffffff434d3f9c-c5e3-4796-b805-00c9b17a6d93.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range("B3"), 'https://github.com/OpenLocalizationTest/oltest/blob/fe32e375100bc6b15da6c00d5bd65ac940e27279/e2e/ffffb6ad0d8b-c1c0-4c1f-bfc6-69772c211997.md', [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, '.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range("D3"), 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c6f150b7c1fd13c10cc737ad4ec2018e6ffe2613/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a3573822-9d94-4592-834d-ffb5b55a6ca7.2924ee7c1e8c42dfa5b48b47664a134ed1cea41e.de-de.xlf', [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, 'a3573822-9d94-4592-834d-ffb5b55a6ca7.2924ee7c1e8c42dfa5b48b47664a134ed1cea41e.de-de.xlf') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F3"), 'https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/505feea415031b0e25c008f05b11f0cfb8e5076b/e2e/a3573822-9d94-4592-834d-ffb5b55a6ca7.md', [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, 'a3573822-9d94-4592-834d-ffb5b55a6ca7.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range("G3"), 'https://github.com/OpenLocalizationTestOrg/olhandback/blob/d20bda3fa226c38ae926feabbe15242408b24720/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a3573822-9d94-4592-834d-ffb5b55a6ca7.2924ee7c1e8c42dfa5b48b47664a134ed1cea41e.de-de.xlf', [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, 'a3573822-9d94-4592-834d-ffb5b55a6ca7.2924ee7c1e8c42dfa5b48b47664a134ed1cea41e.de-de.xlf') | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), 'https://github.com/OpenLocalizationTest/oltest/blob/fe32e375100bc6b15da6c00d5bd65ac940e27279/e2e/ffffff434d3f9c-c5e3-4796-b805-00c9b17a6d93.md', [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, '7a01dfb2-5efe-44c9-9ceb-2d775ef4025f.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range("B4"), 'https://github.com/OpenLocalizationTest/oltest/blob/fe32e375100bc6b15da6c00d5bd65ac940e27279/e2e/ffffff434d3f9c-c5e3-4796-b805-00c9b17a6d93.md', [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, '.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range("D4"), 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c6f150b7c1fd13c10cc737ad4ec2018e6ffe2613/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a3573822-9d94-4592-834d-ffb5b55a6ca7.2924ee7c1e8c42dfa5b48b47664a134ed1cea41e.de-de.xlf', [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, '7a01dfb2-5efe-44c9-9ceb-2d775ef4025f.718e13acabd4baa934d39fc555b8ee9c989bf7e0.de-de.xlf') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F4"), 'https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/505feea415031b0e25c008f05b11f0cfb8e5076b/e2e/a3573822-9d94-4592-834d-ffb5b55a6ca7.md', [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, '7a01dfb2-5efe-44c9-9ceb-2d775ef4025f.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range("G4"), 'https://github.com/OpenLocalizationTestOrg/olhandback/blob/d20bda3fa226c38ae926feabbe15242408b24720/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a3573822-9d94-4592-834d-ffb5b55a6ca7.2924ee7c1e8c42dfa5b48b47664a134ed1cea41e.de-de.xlf', [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, '7a01dfb2-5efe-44c9-9ceb-2d775ef4025f.718e13acabd4baa934d39fc555b8ee9c989bf7e0.de-de.xlf') | Out-Null
